$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '95.918.81'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').Value = "'" + '3.308.20'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('D4').Value = "'" + '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'" + '246.93'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.99%  '
$ws.Range('D6').Value = "'" + '647.28'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.67%  '
$ws.Range('E7').Value = '  -10.45%  '
$ws.Range('E8').Value = '  -5.05%  '
$ws.Range('D9').Value = "'" + '0.999'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = "'" + '0.964'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.69%  '
$ws.Range('D11').Value = "'" + '3.304.96'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.68%  '
$ws.Range('E12').Value = '  -4.68%  '
$ws.Range('D13').Value = "'" + '39.30'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -7.27%  '
$ws.Range('D14').Value = "'" + '95.792.44'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = "'" + '6.03'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -7.88%  '
$ws.Range('E16').Value = '  -4.94%  '
$ws.Range('D17').Value = "'" + '3.927.72'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').Value = "'" + '8.42'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.41%  '
$ws.Range('D19').Value = "'" + '3.321.04'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('D20').Value = "'" + '16.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.15%  '
$ws.Range('B21').Value = 'Stellar'
$ws.Range('C21').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D21').Value = "'" + '0.475'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -8.51%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = "'" + '499.32'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('E23').Value = '  -5.58%  '
$ws.Range('E24').Value = '  -5.56%  '
$ws.Range('E25').Value = '  -6.67%  '
$ws.Range('D26').Value = "'" + '6.33'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('D27').Value = "'" + '93.81'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.42%  '
$ws.Range('D28').Value = "'" + '11.84'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -7.64%  '
$ws.Range('D29').Value = "'" + '3.484.39'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').Value = "'" + '0.140'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -9.05%  '
$ws.Range('D32').Value = "'" + '10.69'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.05%  '
$ws.Range('E33').Value = '  -7.73%  '
$ws.Range('D34').Value = "'" + '2.44'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +10.28%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  -6.65%  '
$ws.Range('D37').Value = "'" + '27.62'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -8.05%  '
$ws.Range('D38').Value = "'" + '1.47'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('D39').Value = "'" + '7.40'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -6.23%  '
$ws.Range('E41').Value = '  -6.10%  '
$ws.Range('D42').Value = "'" + '499.21'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').Value = "'" + '24.28'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.78%  '
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').Value = "'" + '0.0403'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.10%  '
$ws.Range('D47').Value = "'" + '5.35'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('D51').Value = "'" + '3.08'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -7.43%  '
